$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ATECO_code" header (H1) carried a special bold/white-on-theme font
# (used for readability against a colored fill elsewhere in the original
# template). Bring it in line with the rest of the header row's plain
# bold style before renaming, since the cleaned-up header no longer wants
# that distinct formatting.
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Rename the two mislabeled table/header columns to their lowercase,
# snake_case equivalents used by the rest of the schema.
$ws.Range("H1").Value = "ateco_code"
$ws.Range("K1").Value = "cap"

# Restore default (no) selection-range clipboard marquee and move the
# active cell, matching where editing left off.
$excel.CutCopyMode = $false
$ws.Range("G14").Select() | Out-Null
